$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update end time of the second entry (row 3): 17:00 -> 18:30
$ws.Range("E3").Value = 0.77083333333333337

# Insert a new row before the current row 4 (the old blank placeholder row)
# so that the summary rows (sum [min]/[h]/[working weeks]) shift down by one,
# and fill the new row 4 with a third time entry for the same day.
$ws.Rows(4).Insert()

$ws.Range("A4").Value = 2014
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 17
$ws.Range("D4").Value = 0.79166666666666663
$ws.Range("E4").Value = 0.875
$ws.Range("F4").Formula = "=(E4-D4)*24*60"

# Copy the styles from row 3 (data row) onto the new row 4 cells
$ws.Range("D3:G3").Copy() | Out-Null
$ws.Range("D4:G4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the values/formula that PasteSpecial(formats) might have left untouched,
# ensuring the formula text is exactly as expected.
$ws.Range("F4").Formula = "=(E4-D4)*24*60"

# The row that used to be row 4 (blank placeholder with styled D/E/F/G cells)
# is now row 5 after the insert - keep it blank, but make sure styles match.
$ws.Range("D5:G5").Style = "Normal"
$ws.Range("D5").NumberFormat = "hh:mm;@"
$ws.Range("E5").NumberFormat = "hh:mm;@"
$ws.Range("F5").NumberFormat = "0"
$ws.Range("G5").NumberFormat = "hh:mm;@"
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()

# Fix up the SUM formula on the "sum [min]" row (now row 6) to include the new row 5
$ws.Range("F6").Formula = "=SUM(F2:F5)"

# Update the selection to G4, matching the diff
$ws.Range("G4").Select() | Out-Null
